$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry a two-row header (row 1 with units, row 2 with
# "Hiver/Ete/Annee" labels) sitting above the plant data (old rows 3-6).
# Collapse that into a single header row by deleting old row 2 - this
# shifts the old data rows (3-6) up to (2-5) as a structural move, so the
# plant data keeps its original values/number styling untouched.
$ws.Rows(2).Delete()

# Clear what's left of the old header remnants in row 1 before laying out
# the new header.
$ws.Range("A1:K1").Clear()

# ---- New header row ----
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("F1").Font.Name = "Arial"
$ws.Range("F1").Font.Size = 9

$ws.Range("G1").Value = "(MW1)"
$ws.Range("G1").Font.Name = "Arial"
$ws.Range("G1").Font.Size = 9

$ws.Range("H1").Value = "(MW2)"
$ws.Range("H1").Font.Name = "Arial"
$ws.Range("H1").Font.Size = 9

$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Font.Size = 9

$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("J1").Font.Name = "Arial"
$ws.Range("J1").Font.Size = 9

$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("K1").Font.Name = "Arial"
$ws.Range("K1").Font.Size = 9

# ---- Selection matches the post-edit state ----
$ws.Range("A2:K2").Select()

Write-Output "done"
